$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Fill in row 8 of the time-tracking table with the new log entry.
$ws.Range("A8").Value = 45952
$ws.Range("B8").Value = 0.45833333333333331
$ws.Range("C8").Value = 0.70833333333333337
$ws.Range("E8").Value = "Couldn't work on this until now due to busy module. Finised chapter 9, which was about using regular expressions to search and replace text. While I do think this is really helpfull for the Chapters that are coming, it feels like using regexes is quite difficult if you do not do it often. Also reading back what you did seems time consuming"

# The row grows to fit the wrapped description text, same as rows 2-7.
$ws.Rows.Item(8).RowHeight = 57.6

# Move the active selection to C9, matching the saved state of the workbook.
$ws.Range("C9").Select() | Out-Null
